# Updates crypto price/symbol data on sheet1 to match the
# "Updated symbol list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    # Force the cell to stay a text value (avoids Excel silently
    # re-interpreting a numeric-looking string as a float, which would
    # both change the cell type and lose exact textual precision).
    $c.NumberFormat = "@"
    $c.Value = $text
    # Drop back to the workbook's default style so we don't leave a
    # stray "Text" number format attached to the cell.
    $c.Style = "Normal"
}

# --- Column D (Price) value-only refreshes ---
Set-TextValue "D2"  "276.23"
Set-TextValue "D3"  "21.15"
Set-TextValue "D6"  "3.552"
Set-TextValue "D7"  "1.527"
Set-TextValue "D8"  "6.565"
Set-TextValue "D11" "0.08306"
Set-TextValue "D12" "0.03506"
Set-TextValue "D13" "0.03164"
Set-TextValue "D14" "0.09146"
Set-TextValue "D16" "0.001647"
Set-TextValue "D17" "0.04672"
Set-TextValue "D18" "0.006236"
Set-TextValue "D19" "0.006224"
Set-TextValue "D22" "3.719"
Set-TextValue "D25" "0.3291"
Set-TextValue "D28" "0.0002731"
Set-TextValue "D40" "0.04745"

# --- Rows 41-43: coin ranking reshuffled (KickToken/BKEXToken/CEJI cycle) ---
$ws.Range("B41").Value = "CEJI"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D41" "0.005290"
$ws.Range("E41").Value = "40CEJICEJI"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.007016"
$ws.Range("E42").Value = "41KickTokenKICK"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1119"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Remaining column D (Price) value-only refreshes ---
Set-TextValue "D44" "0.01136"
Set-TextValue "D45" "0.00006277"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "D47" "0.7217"
Set-TextValue "D48" "0.001396"
Set-TextValue "D49" "0.00001897"
Set-TextValue "D50" "0.01238"
